$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Общая стоимость" (Total cost) column is being dropped from the
# report header row entirely (not just blanked) - remove its value AND
# its style so the cell disappears from the saved sheet.
$ws.Range("F2").Clear()

# Change F1's formatting: it was part of the (now shorter) merged title
# band. Set the horizontal alignment back to general so only vertical
# centering remains on this column's xf (drop "F2"/title centering).
$ws.Range("F1").HorizontalAlignment = 1

# The merged title cell shrinks from A1:F1 to A1:E1 to line up with the
# narrower header row (columns A-E only).
$ws.Range("A1:F1").UnMerge()
$ws.Range("A1:E1").Merge()

# Reflect the new title band in the sheet's active selection.
$ws.Range("A1:E1").Select()
